$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column Q (year 2020), rows 4-14 (row 3 stays empty / formatting only)
$values = @{
    4  = 2020
    5  = 38.6
    6  = 42.4
    7  = 53.2
    8  = 90.6
    9  = 52.6
    10 = 24.5
    11 = 69.1
    12 = 32.2
    13 = 19.1
    14 = 25.2
}

# Copy the formatting from column P (the previous last year column) into column Q
# for every row of the table (header separator row 3 through the data row 14),
# then write in the new 2020 values.
foreach ($row in 3..14) {
    $src = $ws.Cells.Item($row, 16)   # column P
    $dst = $ws.Cells.Item($row, 17)   # column Q
    $src.Copy()
    $dst.PasteSpecial(-4122)
    if ($values.ContainsKey($row)) {
        $dst.Value2 = $values[$row]
    }
}

$excel.CutCopyMode = 0

# Match the new selection recorded in the saved workbook
$ws.Range("R27").Select() | Out-Null
